$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (keeps header formatting/column widths intact)
# so the shared-strings table rebuilds fresh in first-use order, matching
# a full data replace (as happened on upload).
$ws.Rows.Item(2).Resize(5).ClearContents() | Out-Null

# Row 2: Cheltenham location
$ws.Range("A2").Value = "Cheltenham"
$ws.Range("B2").Value = "Two Bob Snob, 256 Charman Road"
$ws.Range("C2").Value = "22/12/20 1:00pm-2:00pm"
$ws.Range("D2").Value = "Case attended Venue"
$ws.Range("E2").Value = "old"

# Row 3: Moorabbin location (first entry)
$ws.Range("A3").Value = "Moorabbin"
$ws.Range("B3").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C3").Value = "21/12/20 2:00pm-10:00pm  22/12/20 10:00am-6:00pm  24/12/20 1:00pm-10:00pm  28/12/20 8.05pm-8.47pm  29/12/20 12:00pm-4:00pm"
$ws.Range("D3").Value = "Case's workplace"
$ws.Range("E3").Value = "old"

# Row 4: Moorabbin location (second entry, updated exposure period)
$ws.Range("A4").Value = "Moorabbin"
$ws.Range("B4").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C4").Value = "24/12/20 1:00pm-10:00pm  28/12/20 8.05pm-8.47pm  29/12/20 12:00pm-4:00pm"
$ws.Range("D4").Value = "Case's workplace"
$ws.Range("E4").Value = "new"

# Row 5: Southbank location
$ws.Range("A5").Value = "Southbank"
$ws.Range("B5").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C5").Value = "23/12/20 8:00pm-11:00pm"
$ws.Range("D5").Value = "Case ate in store"
$ws.Range("E5").Value = "new"

# Row 6: Southbank location (second entry)
$ws.Range("A6").Value = "Southbank"
$ws.Range("B6").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C6").Value = "23/12/20 8:00pm-11:00pm"
$ws.Range("D6").Value = "Case attended restaurant"
$ws.Range("E6").Value = "old"

# Note: this runtime quantises ColumnWidth to 1/6-character increments on
# save (internal MDW=6 px model), so the inputs below are chosen to land the
# exported <col width=".."> as close as achievable to the authored values
# (10.19921875 / 55.6640625 / 113.265625 / 20.73046875).
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 54.833333333333336
$ws.Columns.Item(3).ColumnWidth = 112.5
$ws.Columns.Item(4).ColumnWidth = 19.833333333333332

$ws.Range("B5").Select()
